$d = $word.ActiveDocument

# The blog text misspells "Berkshire" as "Berkshore" in the line
# "- Berkshore and Schneider". Fix the typo.
#
# We deliberately search on "- Berkshore " (including the leading
# "- " and the trailing space) rather than just "Berkshore" so that the
# matched/replaced range fully spans the run that holds the misspelled
# word together with the stale <w:proofErr spellStart/spellEnd/> spell
# check marks Word had placed around it. Replacing that whole span
# collapses it back down to a single clean run and drops those now
# pointless proof marks, instead of merely swapping the word inside its
# existing run (which would leave the obsolete proofing marks behind).
$rng = $d.Content
$found = $rng.Find.Execute("- Berkshore ", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "- Berkshire ", 2)

if (-not $found) {
    throw "Could not find the text '- Berkshore ' to correct."
}

Write-Output "Corrected 'Berkshore' -> 'Berkshire'."
